# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp values on the Overview, zh-cn and
# de-de sheets to reflect a freshly-generated handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row for 19b1bd9f-... -> Latest HO Xliff Generate Date (col G)
$overview.Range("G4").Value = "2016-09-01 06:51:53"

# zh-cn sheet: row for 19b1bd9f-... -> Correspond Handoff Datetime (col H) and
# Correspond Handback DateTime (col K)
$zhcn.Range("H4").Value = "2016-09-01 06:51:49"
$zhcn.Range("K4").Value = "2016-09-01 06:52:18"

# de-de sheet: row for 19b1bd9f-... -> Correspond Handoff Datetime (col H) and
# Correspond Handback DateTime (col K)
$dede.Range("H4").Value = "2016-09-01 06:51:53"
$dede.Range("K4").Value = "2016-09-01 06:52:26"
